$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): B11 6 -> 9, C11 3 -> 2
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Row 12 (Total): B12 90 -> 135, C12 -9 -> -6, E12 "81/168" -> "129/252"
$ws.Range("B12").Value = 135
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "129/252"
